$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right count 5 -> 4, Wrong count -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right total 105 -> 84, Wrong total -2 -> -4
$ws.Range("B12").Value = 84
$ws.Range("C12").Value = -4

# Update the displayed score string to match new total
$ws.Range("E12").Value = "80 / 112"
